$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.776.93"
$ws.Range("E2").Value = "  -4.03%  "
$ws.Range("D3").Value = "3.133.61"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.96"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.13"
$ws.Range("E6").Value = "  -9.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -8.80%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.135.01"
$ws.Range("E9").Value = "  -5.38%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.116"
$ws.Range("E11").Value = "  -8.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -5.73%  "
$ws.Range("D13").Value = "3.681.01"
$ws.Range("E13").Value = "  -5.47%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "63.879.24"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.78"
$ws.Range("E16").Value = "  -6.79%  "
$ws.Range("D17").Value = "3.137.58"
$ws.Range("E17").Value = "  -5.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000154"
$ws.Range("E18").Value = "  -5.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "409.98"
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  -3.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.18"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.00"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.97"
$ws.Range("E25").Value = "  -3.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.203"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.490"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("E28").Value = "  -12.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.66"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.79"
$ws.Range("E32").Value = "  -6.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.37"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  -5.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.42"
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.10"
$ws.Range("E37").Value = "  -6.89%  "
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("D39").Value = "2.675.43"
$ws.Range("E39").Value = "  -6.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.66"
$ws.Range("E40").Value = "  -7.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.54"
$ws.Range("E42").Value = "  -10.70%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.65"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.699"
$ws.Range("E44").Value = "  -7.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0606"
$ws.Range("E45").Value = "  -7.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.27"
$ws.Range("E46").Value = "  -10.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0256"
$ws.Range("E47").Value = "  -5.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "284.99"
$ws.Range("E48").Value = "  -8.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.88"
$ws.Range("E49").Value = "  -9.57%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0973"
$ws.Range("E51").Value = "  -6.62%  "
